# Generate and save output file after processing
# Insert 3 new columns (history, electives, cs) right before the existing
# "arts" column (R), shifting arts and everything after it 3 columns to
# the right (R->U, S->V, ... AE->AH).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert three blank columns at R:T. Existing R:AE (arts ... act75) shift to U:AH.
$ws.Range("R1:T1").EntireColumn.Insert()

# New header cells for the inserted columns, matching the style of the
# other header cells (bold / centered / bordered). Copy formatting from an
# existing header cell since direct Style object assignment isn't honored.
$ws.Range("Q1").Copy()
$ws.Range("R1:T1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("R1").Value = "general_college_subjects.history"
$ws.Range("S1").Value = "general_college_subjects.electives"
$ws.Range("T1").Value = "general_college_subjects.cs"

# New data cells (row 2) for the inserted columns.
$ws.Range("R2").Value = 0
$ws.Range("S2").Value = 0
$ws.Range("T2").Value = 0

# Normalize casing of the importance/consideration text values.
$ws.Range("D2").Value = "considered"
$ws.Range("E2").Value = "considered"
$ws.Range("F2").Value = "not considered"
$ws.Range("G2").Value = "very important"
$ws.Range("H2").Value = "very important"
$ws.Range("I2").Value = "considered"
$ws.Range("J2").Value = "considered"

$wb.Save()
